# Q3 Update - 2025
# Updates the UN-LES refugee-statistics sheet:
#  1. The "short-url" column (B) changes for every data row.
#  2. A handful of the most recent data rows (120-128) get refreshed
#     refugees / asylum_seekers / ooc figures.
#
# All values in this sheet are stored as text (the workbook was produced
# by a CSV->XLSX converter, not by typing numbers into Excel), so every
# numeric-looking value is written with the cell pre-formatted as Text
# ("@") to stop Excel's COM layer from silently re-typing it as a Number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. short-url column: every data row (2-128) gets the new slug ----
$ws.Range("B2:B128").Value = "Nm5lw4"

# --- 2. refreshed stats for rows 120-128 ------------------------------
# Address -> new value (all numeric-looking, so force Text format first)
$changes = [ordered]@{
    "O120" = "19"
    "N121" = "5"
    "N122" = "262"
    "O122" = "37"
    "T122" = "6"
    "N123" = "80"
    "O123" = "72"
    "N124" = "25"
    "O124" = "14"
    "N125" = "11"
    "O125" = "13"
    "N126" = "16"
    "O126" = "12"
    "O127" = "5"
    "N128" = "8"
    "O128" = "14"
}

foreach ($addr in $changes.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$addr]
}
